$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding percentage values as plain text (inlineStr) must stay text.
$textCells = @("C2","D2","E2","C3","D3","E3","C4","D4","E4","C5","D5","E5","C6","D6","E6","D7","E7","E8","E9","E10","E11")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("C2").Value = "67.23%"
$ws.Range("D2").Value = "95.35%"
$ws.Range("E2").Value = "99.09%"

# Row 3
$ws.Range("B3").Value = 48.9
$ws.Range("C3").Value = "30.68%"
$ws.Range("D3").Value = "86.59%"
$ws.Range("E3").Value = "96.73%"

# Row 4
$ws.Range("B4").Value = 40.7
$ws.Range("C4").Value = "1.06%"
$ws.Range("D4").Value = "9.11%"
$ws.Range("E4").Value = "48.23%"

# Row 5
$ws.Range("C5").Value = "0.79%"
$ws.Range("D5").Value = "6.25%"
$ws.Range("E5").Value = "34.61%"

# Row 6
$ws.Range("C6").Value = "0.24%"
$ws.Range("D6").Value = "2.52%"
$ws.Range("E6").Value = "18.23%"

# Row 7
$ws.Range("D7").Value = "0.11%"
$ws.Range("E7").Value = "1.74%"

# Row 8
$ws.Range("E8").Value = "1.13%"

# Row 9
$ws.Range("E9").Value = "0.23%"

# Row 10
$ws.Range("E10").Value = "0.01%"

# Row 11
$ws.Range("E11").Value = "0.01%"
